$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.953.17"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.876.58"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("D5").Value = "'0.7401"
$ws.Range("E5").Value = "  -4.28%  "
$ws.Range("D6").Value = "'242.26"
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.3156"
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("B9").Value = "Solana"
$ws.Range("C9").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D9").Value = "'24.84"
$ws.Range("E9").Value = "  -3.34%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.07170"
$ws.Range("E10").Value = "  -0.73%  "
$ws.Range("D11").Value = "'0.08434"
$ws.Range("E11").Value = "  -4.95%  "
$ws.Range("D12").Value = "'0.7552"
$ws.Range("E12").Value = "  -2.15%  "
$ws.Range("D13").Value = "'5.414"
$ws.Range("E13").Value = "  -0.31%  "
$ws.Range("D14").Value = "1.873.88"
$ws.Range("E14").Value = "  -12.08%  "
$ws.Range("D15").Value = "'92.84"
$ws.Range("E15").Value = "  -1.72%  "
$ws.Range("D16").Value = "29.946.44"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "'6.098"
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("D18").Value = "'13.64"
$ws.Range("E18").Value = "  -2.18%  "
$ws.Range("D19").Value = "'243.83"
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("D20").Value = "'0.000007838"
$ws.Range("E20").Value = "  -0.29%  "
$ws.Range("D21").Value = "'1.000"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "2.115.18"
$ws.Range("E22").Value = "  -7.31%  "
$ws.Range("D23").Value = "'7.987"
$ws.Range("E23").Value = "  -1.99%  "
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("D25").Value = "'0.1560"
$ws.Range("E25").Value = "  -2.52%  "
$ws.Range("D26").Value = "'9.322"
$ws.Range("E26").Value = "  -2.08%  "
$ws.Range("D27").Value = "'164.76"
$ws.Range("E27").Value = "  +1.44%  "
$ws.Range("D28").Value = "'18.66"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("D29").Value = "'2.041"
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'1.478"
$ws.Range("E30").Value = "  +3.41%  "
$ws.Range("D31").Value = "'4.606"
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("D33").Value = "'4.284"
$ws.Range("E33").Value = "  +3.97%  "
$ws.Range("D34").Value = "'0.05340"
$ws.Range("E34").Value = "  -2.70%  "
$ws.Range("D35").Value = "'1.240"
$ws.Range("E35").Value = "  -0.67%  "
$ws.Range("D36").Value = "'0.7608"
$ws.Range("E36").Value = "  +1.13%  "
$ws.Range("D37").Value = "'1.002"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("D38").Value = "'2.699"
$ws.Range("E38").Value = "  -0.65%  "
$ws.Range("D39").Value = "'0.01958"
$ws.Range("E39").Value = "  -0.10%  "
$ws.Range("D40").Value = "'2.755"
$ws.Range("E40").Value = "  -1.16%  "
$ws.Range("D41").Value = "'0.4490"
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("D42").Value = "1.117.03"
$ws.Range("E42").Value = "  +2.26%  "
$ws.Range("D43").Value = "'6.154"
$ws.Range("E43").Value = "  +2.01%  "
$ws.Range("D44").Value = "'72.67"
$ws.Range("E44").Value = "  -1.43%  "
$ws.Range("D45").Value = "'0.8621"
$ws.Range("E45").Value = "  +0.80%  "
$ws.Range("D46").Value = "'1.002"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").Value = "'103.31"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").Value = "'7.708"
$ws.Range("E48").Value = "  +1.25%  "
$ws.Range("B49").Value = "SynthetixNetwork"
$ws.Range("C49").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D49").Value = "'3.078"
$ws.Range("E49").Value = "  +4.30%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'1.846"
$ws.Range("E50").Value = "  -2.25%  "
$ws.Range("D51").Value = "2.014.18"
$ws.Range("E51").Value = "  -3.60%  "
